$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "245.83"
Set-TextValue "D3" "22.12"
Set-TextValue "D4" "5.363"
Set-TextValue "D5" "0.05868"
Set-TextValue "D6" "3.388"
Set-TextValue "D7" "6.381"
Set-TextValue "D8" "0.8122"
Set-TextValue "D9" "0.9635"
Set-TextValue "D10" "0.1421"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.07366"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03482"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue "D13" "0.03034"
Set-TextValue "D14" "4.468"
Set-TextValue "D15" "0.09386"
Set-TextValue "D16" "0.001600"
Set-TextValue "D17" "0.04838"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D18" "0.006254"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D19" "0.004080"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D20" "0.0009881"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D21" "0.00009704"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D22" "3.687"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D23" "2.198"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D24" "0.01118"
$ws.Range("E24").Value = "23OneONEBestin24h"
Set-TextValue "D25" "0.3253"
Set-TextValue "D27" "0.0002472"
Set-TextValue "D40" "0.03855"
Set-TextValue "D41" "0.006623"
Set-TextValue "D42" "0.1072"
Set-TextValue "D43" "0.003001"
Set-TextValue "D44" "0.005760"
Set-TextValue "D45" "0.00005653"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
Set-TextValue "D48" "0.07648"
